$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.529.47'
$ws.Range("E2").Value = '  +3.18%  '
$ws.Range("D3").Value = '1.841.66'
$ws.Range("E3").Value = '  +2.39%  '
$ws.Range("E4").Value = '  +0.32%  '
$ws.Range("D5").Value = "'231.50"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +3.27%  '
$ws.Range("D6").Value = "'0.621"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +3.21%  '
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("D8").Value = "'43.53"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +10.87%  '
$ws.Range("E9").Value = '  +8.44%  '
$ws.Range("E10").Value = '  +5.18%  '
$ws.Range("E11").Value = '  +2.84%  '
$ws.Range("D12").Value = '2.107.64'
$ws.Range("E12").Value = '  +2.39%  '
$ws.Range("D13").Value = '1.847.86'
$ws.Range("E13").Value = '  +2.65%  '
$ws.Range("E14").Value = '  +3.87%  '
$ws.Range("E15").Value = '  +7.48%  '
$ws.Range("D16").Value = "'4.72"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +8.47%  '
$ws.Range("D17").Value = '35.498.34'
$ws.Range("D18").Value = "'70.43"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +3.73%  '
$ws.Range("E19").Value = '  +5.00%  '
$ws.Range("D20").Value = "'244.67"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +2.65%  '
$ws.Range("D21").Value = "'12.06"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +8.87%  '
$ws.Range("E22").Value = '  +14.25%  '
$ws.Range("E23").Value = '  +0.23%  '
$ws.Range("D24").Value = "'2.23"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +4.11%  '
$ws.Range("D25").Value = "'171.84"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.62%  '
$ws.Range("D26").Value = "'8.02"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +4.76%  '
$ws.Range("D27").Value = "'17.84"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.75%  '
$ws.Range("E28").Value = '  +1.07%  '
$ws.Range("D29").Value = "'1.56"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +27.85%  '
$ws.Range("E30").Value = '  +0.31%  '
$ws.Range("D31").Value = '3.300.59'
$ws.Range("E31").Value = '  +35.84%  '
$ws.Range("D32").Value = "'0.0552"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +7.92%  '
$ws.Range("D33").Value = "'4.10"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +7.31%  '
$ws.Range("D34").Value = "'3.94"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +5.44%  '
$ws.Range("D36").Value = "'94.94"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +16.73%  '
$ws.Range("E37").Value = '  +8.96%  '
$ws.Range("E38").Value = '  +7.92%  '
$ws.Range("D39").Value = '1.349.68'
$ws.Range("E39").Value = '  +3.89%  '
$ws.Range("D40").Value = "'15.49"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +11.49%  '
$ws.Range("E41").Value = '  +7.12%  '
$ws.Range("E42").Value = '  +6.01%  '
$ws.Range("E43").Value = '  +7.86%  '
$ws.Range("E44").Value = '  +4.13%  '
$ws.Range("E45").Value = '  +0.93%  '
$ws.Range("D46").Value = "'2.81"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.78%  '
$ws.Range("E47").Value = '  +9.88%  '
$ws.Range("D48").Value = "'0.0519"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.62%  '
$ws.Range("D49").Value = '2.011.61'
$ws.Range("E49").Value = '  +2.66%  '
$ws.Range("E50").Value = '  +0.27%  '
$ws.Range("D51").Value = "'102.53"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.12%  '
